# Automatische test-sync: 2025-06-19 17:58:30
#
# Appends the newest "Sollicitatie" mail-log entry to the bottom of the
# "Logs" sheet (row 32), widens the conditional-formatting ranges on the
# Categorie/Beantwoord columns to cover the new row, and bumps the
# "Sollicitatie" tally on the "Dashboard" sheet from 2 to 3.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# The sheet currently has data through row 31; the new log entry lands on
# the next row.
$newRow = 32

$logs.Cells.Item($newRow, 1).Value = "Sollicitatie marketingfunctie"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$logs.Cells.Item($newRow, 4).Value = "Sollicitatie"
# Column E (Antwoord) is left blank for this entry.
$logs.Cells.Item($newRow, 6).Value = "2025-06-19 17:58:20"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# Extend the conditional-formatting ranges (Categorie + Beantwoord columns)
# so they keep covering the whole data range, now D2:D32 / G2:G32.
$categorieRange = $logs.Range("D2:D32")
$categorieRules = $categorieRange.FormatConditions
for ($i = 1; $i -le $categorieRules.Count; $i++) {
    $categorieRules.Item($i).ModifyAppliesToRange($categorieRange)
}

$beantwoordRange = $logs.Range("G2:G32")
$beantwoordRules = $beantwoordRange.FormatConditions
for ($i = 1; $i -le $beantwoordRules.Count; $i++) {
    $beantwoordRules.Item($i).ModifyAppliesToRange($beantwoordRange)
}

# Refresh the Dashboard's "Sollicitatie" count (row 6) to reflect the new entry.
$dashboard.Cells.Item(6, 2).Value = 3
